$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 377.9565
$ws.Range("I5").Value = 182
$ws.Range("J5").Value = 528.6923
$ws.Range("K5").Value = 182
$ws.Range("L5").Value = 528.6923
$ws.Range("M5").Value = -67
$ws.Range("N5").Value = -758.6923
$ws.Range("H6").Value = 8474.444
$ws.Range("I6").Value = 8474.444
$ws.Range("K6").Value = 25423.332
$ws.Range("M6").Value = -25311.332
$ws.Range("H9").Value = 923.3333
$ws.Range("I9").Value = 1314.5
$ws.Range("J9").Value = 141
$ws.Range("K9").Value = 1314.5
$ws.Range("L9").Value = 141
$ws.Range("M9").Value = -1145.5
$ws.Range("N9").Value = -479
$ws.Range("H11").Value = 114.8125
$ws.Range("I11").Value = 114.8125
$ws.Range("K11").Value = 114.8125
$ws.Range("M11").Value = 25.1875
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H21").Value = 44777.668
$ws.Range("H23").Value = 44777.668
$ws.Range("H38").Value = 278.33334
$ws.Range("I38").Value = 278.33334
$ws.Range("K38").Value = 835.0000200000001
$ws.Range("M38").Value = -463.0000200000001
$ws.Range("H54").Value = 265000
$ws.Range("I54").Value = 265000
$ws.Range("K54").Value = 265000
$ws.Range("M54").Value = -264514
$ws.Range("H106").Value = 12743.929
$ws.Range("I106").Value = 10368.083
$ws.Range("K106").Value = 10368.083
$ws.Range("M106").Value = -9737.083000000001
$ws.Range("H112").Value = 2499.2593
$ws.Range("J112").Value = 2798.5715
$ws.Range("L112").Value = 8395.7145
$ws.Range("N112").Value = -10611.7145
$ws.Range("H132").Value = 13278.577
$ws.Range("I132").Value = 10782.762
$ws.Range("J132").Value = 23761
$ws.Range("K132").Value = 32348.286
$ws.Range("L132").Value = 71283
$ws.Range("M132").Value = -29818.286
$ws.Range("N132").Value = -76343
$ws.Range("H135").Value = 10880.158
$ws.Range("J135").Value = 10146.3
$ws.Range("L135").Value = 91316.7
$ws.Range("N135").Value = -96386.7
$ws.Range("H138").Value = 5463.732
$ws.Range("J138").Value = 5750.534
$ws.Range("L138").Value = 17251.602
$ws.Range("N138").Value = -27531.602
$ws.Range("H141").Value = 4182.7188
$ws.Range("I141").Value = 3725
$ws.Range("J141").Value = 6166.1665
$ws.Range("K141").Value = 11175
$ws.Range("L141").Value = 18498.4995
$ws.Range("M141").Value = -5995
$ws.Range("N141").Value = -28858.4995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3733.5
$ws.Range("I45").Value = 3840.7693
$ws.Range("J45").Value = 3578.5557
$ws.Range("K45").Value = 3840.7693
$ws.Range("L45").Value = 3578.5557
$ws.Range("M45").Value = -3463.7693
$ws.Range("N45").Value = -4332.5557
$ws.Range("H63").Value = 999
$ws.Range("I63").Value = 999
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 999
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -313
$ws.Range("H66").Value = 999
$ws.Range("I66").Value = 999
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 4995
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -1563
$ws.Range("H102").Value = 17870.428
$ws.Range("I102").Value = 3348.9167
$ws.Range("K102").Value = 3348.9167
$ws.Range("M102").Value = -1726.9167
$ws.Range("H110").Value = 6910.8
$ws.Range("I110").Value = 2460
$ws.Range("K110").Value = 2460
$ws.Range("M110").Value = -415

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1792.125
$ws.Range("I99").Value = 1625.6
$ws.Range("J99").Value = 2069.6667
$ws.Range("K99").Value = 1625.6
$ws.Range("L99").Value = 2069.6667
$ws.Range("M99").Value = -127.5999999999999
$ws.Range("N99").Value = -5065.6667
$ws.Range("H105").Value = 11666.667
$ws.Range("I105").Value = 11666.667
$ws.Range("K105").Value = 11666.667
$ws.Range("M105").Value = -9919.666999999999
$ws.Range("H107").Value = 3339.85
$ws.Range("I107").Value = 3400.7
$ws.Range("K107").Value = 3400.7
$ws.Range("M107").Value = -1480.7
$ws.Range("H134").Value = 6928.685
$ws.Range("I134").Value = 2110.1143
$ws.Range("J134").Value = 15805
$ws.Range("K134").Value = 6330.342900000001
$ws.Range("L134").Value = 47415
$ws.Range("M134").Value = -3795.342900000001
$ws.Range("N134").Value = -52485

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2642.2942
$ws.Range("I22").Value = 2179.8
$ws.Range("J22").Value = 3303
$ws.Range("K22").Value = 2179.8
$ws.Range("L22").Value = 3303
$ws.Range("M22").Value = -1829.8
$ws.Range("N22").Value = -4003
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H132").Value = 5524.921
$ws.Range("I132").Value = 1711.4286
$ws.Range("K132").Value = 5134.2858
$ws.Range("M132").Value = -2604.2858
$ws.Range("H141").Value = 390309.16
$ws.Range("J141").Value = 433766.88
$ws.Range("L141").Value = 433766.88
$ws.Range("N141").Value = -444126.88

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 209.46153
$ws.Range("I33").Value = 192.875
$ws.Range("K33").Value = 1157.25
$ws.Range("M33").Value = -874.25
$ws.Range("H46").Value = 2584.6924
$ws.Range("I46").Value = 344.55554
$ws.Range("K46").Value = 1033.66662
$ws.Range("M46").Value = -942.66662
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H122").Value = 7689632.5
$ws.Range("J122").Value = 1496548.6
$ws.Range("L122").Value = 13468937.4
$ws.Range("N122").Value = -13473837.4
$ws.Range("H131").Value = 1424.93
$ws.Range("J131").Value = 1474.4409
$ws.Range("L131").Value = 4423.322700000001
$ws.Range("N131").Value = -14503.3227
$ws.Range("H141").Value = 10562
$ws.Range("I141").Value = 3166.6667
$ws.Range("K141").Value = 9500.000100000001
$ws.Range("M141").Value = -4320.000100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 19206.2
$ws.Range("J39").Value = 19206.2
$ws.Range("L39").Value = 19206.2
$ws.Range("N39").Value = -20270.2
$ws.Range("H52").Value = 41654.582
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 47985.5
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 47985.5
$ws.Range("M52").Value = -9741
$ws.Range("N52").Value = -48503.5
$ws.Range("H70").Value = 9230.182000000001
$ws.Range("I70").Value = 6969.9
$ws.Range("K70").Value = 6969.9
$ws.Range("M70").Value = -6699.9
$ws.Range("H73").Value = 9230.182000000001
$ws.Range("I73").Value = 6969.9
$ws.Range("K73").Value = 6969.9
$ws.Range("M73").Value = -6033.9
$ws.Range("H113").Value = 103049.89
$ws.Range("I113").Value = 115431.125
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 115431.125
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -113261.125
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 6341.8823
$ws.Range("J122").Value = 10257.714
$ws.Range("L122").Value = 30773.142
$ws.Range("N122").Value = -35673.142

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11002.16
$ws.Range("I22").Value = 10276.333
$ws.Range("K22").Value = 10276.333
$ws.Range("M22").Value = -9981.333000000001
$ws.Range("H27").Value = 11002.16
$ws.Range("I27").Value = 10276.333
$ws.Range("K27").Value = 10276.333
$ws.Range("M27").Value = -10169.333
$ws.Range("H132").Value = 9431.758
$ws.Range("I132").Value = 6456.846
$ws.Range("K132").Value = 19370.538
$ws.Range("M132").Value = -16840.538
$ws.Range("H136").Value = 15418.761
$ws.Range("I136").Value = 12637.038
$ws.Range("K136").Value = 37911.114
$ws.Range("M136").Value = -35361.114

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7361.9546
$ws.Range("I132").Value = 3417.468
$ws.Range("J132").Value = 17119.37
$ws.Range("K132").Value = 10252.404
$ws.Range("L132").Value = 51358.11
$ws.Range("M132").Value = -7722.403999999999
$ws.Range("N132").Value = -56418.11
